$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: new quarter columns AU (31/03/2024) and AV (30/06/2024)
$ws.Range("AU1").Value = "31/03/2024"
$ws.Range("AV1").Value = "30/06/2024"
$ws.Range("AU1").Font.Bold = $true
$ws.Range("AU1").HorizontalAlignment = -4108
$ws.Range("AU1").VerticalAlignment = -4160
$ws.Range("AU1").Borders.LineStyle = 1
$ws.Range("AV1").Font.Bold = $true
$ws.Range("AV1").HorizontalAlignment = -4108
$ws.Range("AV1").VerticalAlignment = -4160
$ws.Range("AV1").Borders.LineStyle = 1

# Data rows: AU and AV values for each line item
$ws.Range("AU2").Value = 6329968.128
$ws.Range("AV2").Value = 6525853.184
$ws.Range("AU3").Value = 757264
$ws.Range("AV3").Value = 1101139.968
$ws.Range("AU4").Value = 194024.992
$ws.Range("AV4").Value = 220180.992
$ws.Range("AU5").Value = 0
$ws.Range("AV5").Value = 0
$ws.Range("AU6").Value = 320766.016
$ws.Range("AV6").Value = 399025.984
$ws.Range("AU7").Value = 0
$ws.Range("AV7").Value = 0
$ws.Range("AU8").Value = 0
$ws.Range("AV8").Value = 0
$ws.Range("AU9").Value = 37758
$ws.Range("AV9").Value = 42142
$ws.Range("AU10").Value = 4067
$ws.Range("AV10").Value = 3939
$ws.Range("AU11").Value = 200648
$ws.Range("AV11").Value = 435852
$ws.Range("AU12").Value = 841201.024
$ws.Range("AV12").Value = 769846.976
$ws.Range("AU13").Value = 0
$ws.Range("AV13").Value = 0
$ws.Range("AU14").Value = 0
$ws.Range("AV14").Value = 0
$ws.Range("AU15").Value = 0
$ws.Range("AV15").Value = 0
$ws.Range("AU16").Value = 359259.008
$ws.Range("AV16").Value = 294560
$ws.Range("AU17").Value = 0
$ws.Range("AV17").Value = 0
$ws.Range("AU18").Value = 0
$ws.Range("AV18").Value = 0
$ws.Range("AU19").Value = 113558
$ws.Range("AV19").Value = 113558
$ws.Range("AU20").Value = 12521
$ws.Range("AV20").Value = 11744
$ws.Range("AU21").Value = 0
$ws.Range("AV21").Value = 0
$ws.Range("AU22").Value = 4705334.784
$ws.Range("AV22").Value = 4628150.784
$ws.Range("AU23").Value = 16237
$ws.Range("AV23").Value = 16066
$ws.Range("AU24").Value = 9931
$ws.Range("AV24").Value = 10649
$ws.Range("AU25").Value = 0
$ws.Range("AV25").Value = 0
$ws.Range("AU26").Value = 6329968.128
$ws.Range("AV26").Value = 6525853.184
$ws.Range("AU27").Value = 452468
$ws.Range("AV27").Value = 577027.968
$ws.Range("AU28").Value = 12920
$ws.Range("AV28").Value = 15713
$ws.Range("AU29").Value = 47595
$ws.Range("AV29").Value = 58231
$ws.Range("AU30").Value = 28942
$ws.Range("AV30").Value = 41658
$ws.Range("AU31").Value = 196826
$ws.Range("AV31").Value = 277792.992
$ws.Range("AU32").Value = 0
$ws.Range("AV32").Value = 0
$ws.Range("AU33").Value = 0
$ws.Range("AV33").Value = 0
$ws.Range("AU34").Value = 166184.992
$ws.Range("AV34").Value = 183632.992
$ws.Range("AU35").Value = 0
$ws.Range("AV35").Value = 0
$ws.Range("AU36").Value = 0
$ws.Range("AV36").Value = 0
$ws.Range("AU37").Value = 2112070.016
$ws.Range("AV37").Value = 2162141.952
$ws.Range("AU38").Value = 1766407.04
$ws.Range("AV38").Value = 1775544.064
$ws.Range("AU39").Value = 0
$ws.Range("AV39").Value = 0
$ws.Range("AU40").Value = 196140
$ws.Range("AV40").Value = 248800.992
$ws.Range("AU41").Value = 149523.008
$ws.Range("AV41").Value = 137796.992
$ws.Range("AU42").Value = 0
$ws.Range("AV42").Value = 0
$ws.Range("AU43").Value = 0
$ws.Range("AV43").Value = 0
$ws.Range("AU44").Value = 0
$ws.Range("AV44").Value = 0
$ws.Range("AU45").Value = 0
$ws.Range("AV45").Value = 0
$ws.Range("AU46").Value = 11097
$ws.Range("AV46").Value = 10701
$ws.Range("AU47").Value = 3754333.016
$ws.Range("AV47").Value = 3775981.88
$ws.Range("AU48").Value = 2735382.016
$ws.Range("AV48").Value = 2735382.016
$ws.Range("AU49").Value = -26991
$ws.Range("AV49").Value = 14820
$ws.Range("AU50").Value = 0
$ws.Range("AV50").Value = 0
$ws.Range("AU51").Value = 990787.008
$ws.Range("AV51").Value = 878750.0159999999
$ws.Range("AU52").Value = 55155
$ws.Range("AV52").Value = 147030
$ws.Range("AU53").Value = 0
$ws.Range("AV53").Value = 0
$ws.Range("AU54").Value = 0
$ws.Range("AV54").Value = 0
$ws.Range("AU55").Value = 0
$ws.Range("AV55").Value = 0
$ws.Range("AU56").Value = 0
$ws.Range("AV56").Value = 0
$ws.Range("AU59").Value = 53841
$ws.Range("AV59").Value = 53401
$ws.Range("AU60").Value = -1453
$ws.Range("AV60").Value = -1296
$ws.Range("AU61").Value = 52388
$ws.Range("AV61").Value = 52105
$ws.Range("AU62").Value = -2548
$ws.Range("AV62").Value = -2566
$ws.Range("AU63").Value = -12079
$ws.Range("AV63").Value = -11072
$ws.Range("AU64").Value = 0
$ws.Range("AV64").Value = 0
$ws.Range("AU65").Value = 35797
$ws.Range("AV65").Value = 350468.992
$ws.Range("AU66").Value = -2520
$ws.Range("AV66").Value = -249682
$ws.Range("AU67").Value = 1823
$ws.Range("AV67").Value = 271
$ws.Range("AU68").Value = -9123
$ws.Range("AV68").Value = -32455
$ws.Range("AU69").Value = 30499
$ws.Range("AV69").Value = 29016
$ws.Range("AU70").Value = -39622
$ws.Range("AV70").Value = -61471
$ws.Range("AU74").Value = 63738
$ws.Range("AV74").Value = 107070
$ws.Range("AU75").Value = -4810
$ws.Range("AV75").Value = -17304
$ws.Range("AU76").Value = -3628
$ws.Range("AV76").Value = 2120
$ws.Range("AU79").Value = -145
$ws.Range("AV79").Value = -13
$ws.Range("AU80").Value = 55155
$ws.Range("AV80").Value = 91873
